# Auto-generated edit script: updates market-price data cells (columns H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR "Leve Profits" tables, as produced
# by the scheduled data-refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 477119.53
$ws.Range("I28").Value = 769844.0600000001
$ws.Range("J28").Value = 1442.125
$ws.Range("K28").Value = 769844.0600000001
$ws.Range("L28").Value = 1442.125
$ws.Range("M28").Value = -769359.0600000001
$ws.Range("N28").Value = -2412.125

$ws.Range("H43").Value = 13029.1875
$ws.Range("I43").Value = 7499.6665
$ws.Range("J43").Value = 16346.9
$ws.Range("K43").Value = 7499.6665
$ws.Range("L43").Value = 16346.9
$ws.Range("M43").Value = -7430.6665
$ws.Range("N43").Value = -16484.9

$ws.Range("H55").Value = 306.7143
$ws.Range("I55").Value = 79
$ws.Range("J55").Value = 876
$ws.Range("K55").Value = 79
$ws.Range("L55").Value = 876
$ws.Range("M55").Value = 135
$ws.Range("N55").Value = -1304

$ws.Range("H112").Value = 2788.2
$ws.Range("J112").Value = 2788.2
$ws.Range("L112").Value = 8364.599999999999
$ws.Range("N112").Value = -10580.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 10000
$ws.Range("J3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("N3").Value = -10230

$ws.Range("H4").Value = 388.6
$ws.Range("I4").Value = 313.76923
$ws.Range("J4").Value = 875
$ws.Range("K4").Value = 313.76923
$ws.Range("L4").Value = 875
$ws.Range("M4").Value = -197.76923
$ws.Range("N4").Value = -1107

$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("L8").ClearContents()

$ws.Range("H14").Value = 128999.6
$ws.Range("I14").Value = 201666.33
$ws.Range("K14").Value = 201666.33
$ws.Range("M14").Value = -201491.33

$ws.Range("H16").Value = 24942.5
$ws.Range("I16").Value = 24942.5
$ws.Range("K16").Value = 24942.5
$ws.Range("M16").Value = -24655.5

$ws.Range("H21").Value = 4000
$ws.Range("I21").Value = 4000
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 4000
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = -3626
$ws.Range("M21").ClearContents()

$ws.Range("H26").Value = 13246
$ws.Range("I26").Value = 13993
$ws.Range("K26").Value = 13993
$ws.Range("M26").Value = -13663

$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("L33").ClearContents()

$ws.Range("H61").Value = 3255.5
$ws.Range("I61").Value = 2572.3
$ws.Range("J61").Value = 3635.0557
$ws.Range("K61").Value = 2572.3
$ws.Range("L61").Value = 3635.0557
$ws.Range("M61").Value = -2360.3
$ws.Range("N61").Value = -4059.0557

$ws.Range("H74").Value = 3455.0256
$ws.Range("I74").Value = 2883.5
$ws.Range("J74").Value = 3652.1035
$ws.Range("K74").Value = 2883.5
$ws.Range("L74").Value = 3652.1035
$ws.Range("M74").Value = -2009.5
$ws.Range("N74").Value = -5400.1035

$ws.Range("H76").Value = 30532.777
$ws.Range("J76").Value = 30532.777
$ws.Range("L76").Value = 30532.777
$ws.Range("N76").Value = -31208.777

$ws.Range("H77").Value = 3455.0256
$ws.Range("I77").Value = 2883.5
$ws.Range("J77").Value = 3652.1035
$ws.Range("K77").Value = 14417.5
$ws.Range("L77").Value = 18260.5175
$ws.Range("M77").Value = -10049.5
$ws.Range("N77").Value = -26996.5175

$ws.Range("H79").Value = 30532.777
$ws.Range("J79").Value = 30532.777
$ws.Range("L79").Value = 30532.777
$ws.Range("N79").Value = -32872.777

$ws.Range("H132").Value = 4481.4185
$ws.Range("I132").Value = 4823.2354
$ws.Range("K132").Value = 14469.7062
$ws.Range("M132").Value = -11939.7062

$ws.Range("H136").Value = 3255.5
$ws.Range("I136").Value = 2572.3
$ws.Range("J136").Value = 3635.0557
$ws.Range("K136").Value = 7716.900000000001
$ws.Range("L136").Value = 10905.1671
$ws.Range("M136").Value = -5166.900000000001
$ws.Range("N136").Value = -16005.1671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 60607976
$ws.Range("I99").Value = 74075580
$ws.Range("J99").Value = 3722
$ws.Range("K99").Value = 74075580
$ws.Range("L99").Value = 3722
$ws.Range("M99").Value = -74074082
$ws.Range("N99").Value = -6718

$ws.Range("H134").Value = 3318.9666
$ws.Range("I134").Value = 3224.9583
$ws.Range("J134").Value = 3695
$ws.Range("K134").Value = 9674.874899999999
$ws.Range("L134").Value = 11085
$ws.Range("M134").Value = -7139.874899999999
$ws.Range("N134").Value = -16155

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2300
$ws.Range("I134").Value = 1950
$ws.Range("K134").Value = 5850
$ws.Range("M134").Value = -3315

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 5430222.5
$ws.Range("I23").Value = 388.83334
$ws.Range("J23").Value = 8688123
$ws.Range("K23").Value = 1166.50002
$ws.Range("L23").Value = 26064369
$ws.Range("M23").Value = -931.5000199999999
$ws.Range("N23").Value = -26064839

$ws.Range("H68").Value = 1999
$ws.Range("I68").Value = 2002
$ws.Range("J68").Value = 1998.25
$ws.Range("K68").Value = 6006
$ws.Range("L68").Value = 5994.75
$ws.Range("M68").Value = -5195
$ws.Range("N68").Value = -7616.75

$ws.Range("H71").Value = 1999
$ws.Range("I71").Value = 2002
$ws.Range("J71").Value = 1998.25
$ws.Range("K71").Value = 18018
$ws.Range("L71").Value = 17984.25
$ws.Range("M71").Value = -13962
$ws.Range("N71").Value = -26096.25

$ws.Range("H97").Value = 1953.25
$ws.Range("J97").Value = 1562.7142
$ws.Range("L97").Value = 4688.142599999999
$ws.Range("N97").Value = -5680.142599999999

$ws.Range("H107").Value = 962.1818
$ws.Range("I107").Value = 633
$ws.Range("J107").Value = 1035.3334
$ws.Range("K107").Value = 1899
$ws.Range("L107").Value = 3106.0002
$ws.Range("M107").Value = 21
$ws.Range("N107").Value = -6946.0002

$ws.Range("H139").Value = 4304.4
$ws.Range("I139").Value = 3380.5
$ws.Range("J139").Value = 8000
$ws.Range("K139").Value = 10141.5
$ws.Range("L139").Value = 24000
$ws.Range("M139").Value = -5001.5
$ws.Range("N139").Value = -34280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 120.318184
$ws.Range("J2").Value = 73
$ws.Range("L2").Value = 73
$ws.Range("N2").Value = -299

$ws.Range("H3").Value = 9000000
$ws.Range("J3").Value = 9000000
$ws.Range("L3").Value = 9000000
$ws.Range("N3").Value = -9000232

$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("L25").ClearContents()
$ws.Range("M25").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("N50").Value = 0
$ws.Range("L50").ClearContents()
$ws.Range("M50").ClearContents()

$ws.Range("H132").Value = 5020.24
$ws.Range("I132").Value = 4955.7144
$ws.Range("J132").Value = 5170.8
$ws.Range("K132").Value = 14867.1432
$ws.Range("L132").Value = 15512.4
$ws.Range("M132").Value = -12337.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 11497.5
$ws.Range("I2").Value = 11497.5
$ws.Range("K2").Value = 11497.5
$ws.Range("M2").Value = -11385.5

$ws.Range("H3").Value = 4048.25
$ws.Range("I3").Value = 4999.5
$ws.Range("J3").Value = 3097
$ws.Range("K3").Value = 4999.5
$ws.Range("L3").Value = 3097
$ws.Range("M3").Value = -4885.5
$ws.Range("N3").Value = -3325

$ws.Range("H8").Value = 1666.6666
$ws.Range("I8").Value = 1750
$ws.Range("J8").Value = 1500
$ws.Range("K8").Value = 1750
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = -1610
$ws.Range("N8").Value = -1780

$ws.Range("H126").Value = 32260396
$ws.Range("I126").Value = 2169.875
$ws.Range("K126").Value = 6509.625
$ws.Range("M126").Value = -4039.625

$ws.Range("H132").Value = 1289.5555
$ws.Range("I132").Value = 738.9524
$ws.Range("J132").Value = 3216.6667
$ws.Range("K132").Value = 2216.8572
$ws.Range("L132").Value = 9650.000100000001
$ws.Range("M132").Value = 313.1428000000001
$ws.Range("N132").Value = -14710.0001
